$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Segundo Caso de prueba fallo
# 1. Fill in the second "Resultado esperado" row (B22, merged B22:H22) with the
#    failure-scenario description.
$ws.Range("B22").Value = "El candado no se abre a pesar de que se ingresó la combinación correcta o se utilizó la llave adecuada."

# 2. Update the "Estado" cell (B24) from "Exitoso" to "fallo", keeping the bold
#    "Estado:" label intact and only replacing the status word itself.
$cell = $ws.Range("B24")
$cell.Characters(9, 7).Text = "fallo"
$cell.Characters(1, 7).Font.Bold = $true
$cell.Characters(1, 7).Font.Size = 12
$cell.Characters(8, 9).Font.Size = 11

# 3. Widen column H a bit (cosmetic change made alongside the content edit).
$ws.Columns("H").ColumnWidth = 21.6

# 4. Leave the selection where the author left it when saving.
$ws.Range("D23").Select()
